# Updates cryptos list figures (price and 1h volume change) per the
# latest scrape, matching the authored commit "Updated cryptos list on
# Thu Aug  3 22:53:33 UTC 2023 with GitHub Actions".
#
# Column D ("Price") cells are plain text in the source data (values such
# as "29.198.30" or "1.835.22" are not valid numbers, and even the
# numeric-looking ones like "22.63" must stay as text to match the
# original formatting), so we force the NumberFormat to text ("@")
# before writing each Price value to stop Excel from re-interpreting it
# as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.198.30'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.835.22'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.76'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6660'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07364'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2919'
$ws.Range('E9').Value = '  -2.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.63'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07697'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.823.88'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.973'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6653'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.48'
$ws.Range('E15').Value = '  -4.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.085'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.132.43'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008261'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '225.86'
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.43'
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.116'
$ws.Range('E22').Value = '  -3.76%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.62'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.633'
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1391'
$ws.Range('E26').Value = '  -4.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.91'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.503'
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.109'
$ws.Range('E29').Value = '  -3.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.026'
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05294'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.868'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7531'
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.130'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.301.28'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01796'
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.718'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9205'
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('B41').Value = 'XinFinNetwork'
$ws.Range('C41').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.08564'
$ws.Range('E41').Value = '  +16.03%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.942'
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.26'
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000127'
$ws.Range('E45').Value = '  +3.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.970.71'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '63.28'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05931'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.995'
$ws.Range('E51').Value = '  -5.65%  '
